$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.73"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.47"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.72%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.076"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.35%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06660"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.388"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.77%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.403"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.00%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.374"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7.43%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9372"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.22%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1576"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.60%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06471"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.18%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07640"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.29%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02937"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.14%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.43%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001575"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.68%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04500"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.06%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006447"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.24%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006257"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.70%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.63%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.250"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.97%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3216"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.25%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.87%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.071"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.09%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1551"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.84%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001191"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.08%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004139"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.62%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001248"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.56%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001616"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.85%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04218"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.50%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006722"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.08%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.52%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001967"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.18%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01236"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.93%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005556"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.52%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.966"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "25.93%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01305"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-28.94%"
